$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 10 & 11: swap the A (Id), Q (Ost) and R (Nord) values between them ---
$swapCols1011 = @("A", "Q", "R")
foreach ($col in $swapCols1011) {
    $addr10 = "$col" + "10"
    $addr11 = "$col" + "11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2
    $ws.Range($addr10).Value2 = $v11
    $ws.Range($addr11).Value2 = $v10
}

# --- Rows 12 & 13: swap the record contents between them (only the columns   ---
# --- that actually differ between the two records: id/taxon/location/AX)    ---
$swapCols1213 = @("A", "B", "D", "E", "F", "G", "H", "K", "L", "M", "N", "Q", "R", "AX")
foreach ($col in $swapCols1213) {
    $addr12 = "$col" + "12"
    $addr13 = "$col" + "13"
    $v12 = $ws.Range($addr12).Value2
    $v13 = $ws.Range($addr13).Value2
    $ws.Range($addr12).Value2 = $v13
    $ws.Range($addr13).Value2 = $v12
}
